{"js": "// The placeholder run \"${schule}\" must become \"${schule_nametype}\",\n// ending up split across three runs (same run formatting throughout):\n//   \"${schule\"  +  \"_nametype\"  +  \"}\"\n// This mirrors how Word itself splits a run when a formatting change is\n// (re-)applied to only part of it.\n\nconst body = context.document.body;\n\n// 1) Locate the \"${schule}\" placeholder and replace its text in place.\n//    insertText(..., \"Replace\") keeps the run's existing formatting and\n//    preserves the surrounding bookmark (_GoBack) wrapping it.\nconst placeholder = body.search(\"${schule}\", { matchCase: true, matchWholeWord: false });\nplaceholder.load(\"items\");\nawait context.sync();\n\nif (placeholder.items.length === 0) {\n  throw new Error('Could not find \"${schule}\" placeholder in the document body.');\n}\n\nplaceholder.items[0].insertText(\"${schule_nametype}\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Re-locate the newly inserted \"_nametype\" segment and toggle its bold\n//    formatting off then back on. Re-applying identical run formatting to\n//    only part of the run forces Word to split it into separate runs at\n//    that boundary, which is exactly the run layout the target XML has:\n//    \"${schule\" / \"_nametype\" / \"}\" as three sibling runs.\nconst nametype = body.search(\"_nametype\", { matchCase: true, matchWholeWord: false });\nnametype.load(\"items\");\nawait context.sync();\n\nif (nametype.items.length === 0) {\n  throw new Error('Could not find \"_nametype\" segment after inserting it.');\n}\n\nconst nametypeRange = nametype.items[0];\nnametypeRange.font.load(\"bold\");\nawait context.sync();\n\nnametypeRange.font.bold = false;\nawait context.sync();\n\nnametypeRange.font.bold = true;\nawait context.sync();\n", "ps1": "# The placeholder run \"${schule}\" must become \"${schule_nametype}\",\n# ending up split across three runs (same run formatting throughout):\n#   \"${schule\"  +  \"_nametype\"  +  \"}\"\n# This mirrors how Word itself splits a run when a formatting change is\n# (re-)applied to only part of it.\n\n$d = $word.ActiveDocument\n\n# 1) Locate the \"${schule}\" placeholder and replace its text in place.\n#    Assigning Range.Text keeps the run's existing formatting and\n#    preserves the surrounding bookmark (_GoBack) wrapping it.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"`${schule}\"\n$find.MatchWildcards = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"${schule}\" placeholder in the document.'\n}\n\n$placeholder = $find.Parent\n$placeholder.Text = \"`${schule_nametype}\"\n\n# 2) Re-locate the newly inserted \"_nametype\" segment and toggle its bold\n#    formatting off then back on. Re-applying identical run formatting to\n#    only part of the run forces Word to split it into separate runs at\n#    that boundary, which is exactly the run layout the target XML has:\n#    \"${schule\" / \"_nametype\" / \"}\" as three sibling runs.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"_nametype\"\n$find2.MatchWildcards = $false\n$found2 = $find2.Execute()\n\nif (-not $found2) {\n    throw 'Could not find \"_nametype\" segment after inserting it.'\n}\n\n$nametypeRange = $find2.Parent\n$nametypeRange.Font.Bold = 0\n$nametypeRange.Font.Bold = 1\n"}
